$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows to append (Transaction Type, Payment Type, InternalComment, USD Amount)
$newRows = @(
    @{ Row = 30; E = "Withdrawal"; N = "Credit Card"; P = "Tradeprof"; T = 271.6875 },
    @{ Row = 31; E = "Deposit";    N = "Wiretransfer"; P = "Roobic";    T = 3180.0428000000002 },
    @{ Row = 32; E = "Deposit";    N = "Crypto";       P = "ETH";       T = 4061.92 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E    # column E
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # column N
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # column P
    $ws.Cells.Item($r.Row, 20).Value = $r.T   # column T
}

# Update selection to match the new viewport
$ws.Range("Q27:S36").Select()
